# Auto-generated PowerShell Word COM-interop script
$d = $word.ActiveDocument

# 1. Remove the old _GoBack bookmark from "List append = O(1)" paragraph
$d.Bookmarks.ShowHidden = $true
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Convert the trailing YouTube URL paragraph into a real hyperlink
$lastParaIndex = $d.Paragraphs.Count
$urlPara = $d.Paragraphs($lastParaIndex)
$urlRange = $d.Range($urlPara.Range.Start, $urlPara.Range.End - 1)
$url = $urlRange.Text
$d.Hyperlinks.Add($urlRange, $url) | Out-Null

# 3. Append all of the new paragraphs
$cur = $d.Paragraphs($d.Paragraphs.Count)
$cur.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "BST search trees have some cool traits"
$p.Style = "Heading1"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Versus a heap"
$p.Style = "Heading2"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "A heap is an array that can be visualized as a tree"
$p.Style = "Heading3"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "But a tree actually has pointers"
$p.Style = "Heading3"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Pointers"
$p.Style = "Heading2"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Parent"
$p.Style = "Heading3"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Left"
$p.Style = "Heading3"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Right"
$p.Style = "Heading3"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "If you’re in the left subtree"
$p.Style = "Heading2"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Your values are less than those of the right subtree"
$p.Style = "Heading3"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "And vice versa"
$p.Style = "Heading3"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "It’s essentially a sorted list + sorted array"
$p.Style = "Heading2"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "During the insertion phase"
$p.Style = "Heading2"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "You can choose to make checks"
$p.Style = "Heading3"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "THIS IS THE MOST FLEXIBLE THING OF A BST"
$p.Style = "Heading3"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Insertion"
$p.Style = "Heading1"
$bmRange = $d.Range($p.Range.Start, $p.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "O(h) where h is the height of the tree"
$p.Style = "Heading2"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Find Min"
$p.Style = "Heading1"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Keep going to the left"
$p.Style = "Heading2"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "O(h)"
$p.Style = "Heading2"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Find Max"
$p.Style = "Heading1"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Keep going to the right"
$p.Style = "Heading2"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "O(h)"
$p.Style = "Heading2"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Augmented Binary Search Trees"
$p.Style = "Heading1"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "You can put more than one value in"
$p.Style = "Heading2"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "I.e. you can add how many children are underneath"
$p.Style = "Heading2"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Toy Problem, find all values <= x"
$p.Style = "Heading1"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Assuming you have an augmented BST with size associated with it"
$p.Style = "Heading2"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Every time you look arrive at a new node"
$p.Style = "Heading2"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "+1"
$p.Style = "Heading3"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "If you see something is higher and you go to it"
$p.Style = "Heading2"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "+children in left node"
$p.Style = "Heading3"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Will cover everything in remaining left branch"
$p.Style = "Heading4"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "When you finally arrive at your node, do one last check"
$p.Style = "Heading2"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "+children in left node as necessary"
$p.Style = "Heading3"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Or keep going until you reach double null while still being less than your desired endpoint"
$p.Style = "Heading2"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "|*.*||*.*|BFS|*.*||*.*|"
$p.Style = "Heading1"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "|*.*||*.*|DFS|*.*||*.*|"
$p.Style = "Heading1"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "HARD PROBLEM"
$p.Style = "Heading2"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "https://leetcode.com/problems/critical-connections-in-a-network/discuss/494896/Python-DFS"
$p.Style = "Heading3"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Edges"
$p.Style = "Heading2"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "In an Undirected graph"
$p.Style = "Heading3"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Each edge is visited twice"
$p.Style = "Heading4"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "In a Directed graph"
$p.Style = "Heading3"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Each edge visited  once"
$p.Style = "Heading4"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Tree edge"
$p.Style = "Heading3"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "When you visit a NEW vertex via an edge"
$p.Style = "Heading4"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Because they form a directed tree"
$p.Style = "Heading4"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Back edge"
$p.Style = "Heading3"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Directed or reaches an ancestor"
$p.Style = "Heading4"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "Back edge"
$p.Style = "Heading2"
$p.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count

$p = $d.Paragraphs($idx)
$p.Range.Text = "There’s a tree edge"
$p.Style = "Heading3"
